# Auto-generated edit script to update Excalibur_Profits sheets per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1047.5
$ws.Range("I131").Value = 1047.5
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3142.5
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 1897.5
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 228350.53
$ws.Range("I132").Value = 243495.45
$ws.Range("J132").Value = 4205.8
$ws.Range("K132").Value = 730486.3500000001
$ws.Range("L132").Value = 12617.4
$ws.Range("M132").Value = -727956.3500000001
$ws.Range("N132").Value = -17677.4
$ws.Range("H137").Value = 2093.549
$ws.Range("I137").Value = 1110.591
$ws.Range("J137").Value = 2839.2415
$ws.Range("K137").Value = 3331.773
$ws.Range("L137").Value = 8517.7245
$ws.Range("M137").Value = -781.7729999999997
$ws.Range("N137").Value = -13617.7245
$ws.Range("H138").Value = 3370.843
$ws.Range("I138").Value = 2562.4285
$ws.Range("J138").Value = 3676.7297
$ws.Range("K138").Value = 7687.2855
$ws.Range("L138").Value = 11030.1891
$ws.Range("M138").Value = -2547.2855
$ws.Range("N138").Value = -21310.1891

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1301.8
$ws.Range("I25").Value = 1069.6666
$ws.Range("K25").Value = 1069.6666
$ws.Range("M25").Value = -667.6666
$ws.Range("H110").Value = 2422.75
$ws.Range("I110").Value = 1976.7222
$ws.Range("K110").Value = 1976.7222
$ws.Range("M110").Value = 68.27780000000007
$ws.Range("H132").Value = 3909.2644
$ws.Range("I132").Value = 2460.8147
$ws.Range("J132").Value = 6279.4546
$ws.Range("K132").Value = 7382.4441
$ws.Range("L132").Value = 18838.3638
$ws.Range("M132").Value = -4852.4441
$ws.Range("N132").Value = -23898.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1365.0667
$ws.Range("I122").Value = 1405.5
$ws.Range("K122").Value = 4216.5
$ws.Range("M122").Value = -1766.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 151540.22
$ws.Range("I2").Value = 200023.36
$ws.Range("K2").Value = 1200140.16
$ws.Range("M2").Value = -1200027.16
$ws.Range("H6").Value = 128.625
$ws.Range("I6").Value = 113.166664
$ws.Range("J6").Value = 175
$ws.Range("K6").Value = 339.499992
$ws.Range("L6").Value = 525
$ws.Range("M6").Value = -226.499992
$ws.Range("N6").Value = -751
$ws.Range("H7").Value = 5730056
$ws.Range("I7").Value = 10000094
$ws.Range("J7").Value = 4022040.5
$ws.Range("K7").Value = 30000282
$ws.Range("L7").Value = 12066121.5
$ws.Range("M7").Value = -30000170
$ws.Range("N7").Value = -12066345.5
$ws.Range("H10").Value = 1081.3334
$ws.Range("I10").Value = 717.6
$ws.Range("K10").Value = 2152.8
$ws.Range("M10").Value = -2013.8
$ws.Range("H11").Value = 924.1429000000001
$ws.Range("I11").Value = 486.4
$ws.Range("K11").Value = 1459.2
$ws.Range("M11").Value = -1319.2
$ws.Range("H12").Value = 20
$ws.Range("J12").Value = 20
$ws.Range("L12").Value = 60
$ws.Range("N12").Value = -406
$ws.Range("H13").Value = 1549.6666
$ws.Range("J13").Value = 1849.6
$ws.Range("L13").Value = 5548.799999999999
$ws.Range("N13").Value = -5884.799999999999
$ws.Range("H15").Value = 177.6
$ws.Range("J15").Value = 260.66666
$ws.Range("L15").Value = 781.9999799999999
$ws.Range("N15").Value = -1061.99998
$ws.Range("H16").Value = 1082.875
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 1706.2
$ws.Range("K16").Value = 132
$ws.Range("L16").Value = 5118.6
$ws.Range("M16").Value = 41
$ws.Range("N16").Value = -5464.6
$ws.Range("H33").Value = 715.8333
$ws.Range("I33").Value = 566.6667
$ws.Range("K33").Value = 3400.0002
$ws.Range("M33").Value = -3117.0002
$ws.Range("H44").Value = 364.8
$ws.Range("I44").Value = 175
$ws.Range("J44").Value = 491.33334
$ws.Range("K44").Value = 525
$ws.Range("L44").Value = 1474.00002
$ws.Range("M44").Value = -127
$ws.Range("N44").Value = -2270.00002
$ws.Range("H70").Value = 3676.75
$ws.Range("I70").Value = 2237.3333
$ws.Range("K70").Value = 6711.999899999999
$ws.Range("M70").Value = -6396.999899999999
$ws.Range("H73").Value = 3676.75
$ws.Range("I73").Value = 2237.3333
$ws.Range("K73").Value = 6711.999899999999
$ws.Range("M73").Value = -5619.999899999999
$ws.Range("H75").Value = 3970.9333
$ws.Range("J75").Value = 4866.4546
$ws.Range("L75").Value = 14599.3638
$ws.Range("N75").Value = -16595.3638
$ws.Range("H78").Value = 3970.9333
$ws.Range("J78").Value = 4866.4546
$ws.Range("L78").Value = 43798.0914
$ws.Range("N78").Value = -53782.0914
$ws.Range("H80").Value = 5473.25
$ws.Range("J80").Value = 5632.6665
$ws.Range("L80").Value = 16897.9995
$ws.Range("N80").Value = -18769.9995
$ws.Range("H81").Value = 90913770
$ws.Range("J81").Value = 6319.875
$ws.Range("L81").Value = 18959.625
$ws.Range("N81").Value = -21205.625
$ws.Range("H82").Value = 8982
$ws.Range("H83").Value = 5473.25
$ws.Range("J83").Value = 5632.6665
$ws.Range("L83").Value = 50693.9985
$ws.Range("N83").Value = -60053.9985
$ws.Range("H84").Value = 90913770
$ws.Range("J84").Value = 6319.875
$ws.Range("L84").Value = 56878.875
$ws.Range("N84").Value = -68110.875
$ws.Range("H85").Value = 8982
$ws.Range("H86").Value = 1127.5454
$ws.Range("J86").Value = 1250.5
$ws.Range("L86").Value = 3751.5
$ws.Range("N86").Value = -6123.5
$ws.Range("H87").Value = 23012.8
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 23012.8
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 69038.39999999999
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -71534.39999999999
$ws.Range("H88").Value = 3998.5
$ws.Range("J88").Value = 3998.5
$ws.Range("L88").Value = 11995.5
$ws.Range("N88").Value = -12851.5
$ws.Range("H89").Value = 1127.5454
$ws.Range("J89").Value = 1250.5
$ws.Range("L89").Value = 11254.5
$ws.Range("N89").Value = -23110.5
$ws.Range("H90").Value = 23012.8
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 23012.8
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 207115.2
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -219595.2
$ws.Range("H91").Value = 3998.5
$ws.Range("J91").Value = 3998.5
$ws.Range("L91").Value = 11995.5
$ws.Range("N91").Value = -14959.5
$ws.Range("H107").Value = 498.15384
$ws.Range("I107").Value = 456.25
$ws.Range("J107").Value = 516.7778
$ws.Range("K107").Value = 1368.75
$ws.Range("L107").Value = 1550.3334
$ws.Range("M107").Value = 551.25
$ws.Range("N107").Value = -5390.3334
$ws.Range("H113").Value = 2068.6897
$ws.Range("J113").Value = 2095.8462
$ws.Range("L113").Value = 6287.5386
$ws.Range("N113").Value = -10627.5386
$ws.Range("H121").Value = 824.75
$ws.Range("J121").Value = 844.8889
$ws.Range("L121").Value = 2534.6667
$ws.Range("N121").Value = -5154.6667
$ws.Range("H131").Value = 12531.15
$ws.Range("I131").Value = 777.5714
$ws.Range("J131").Value = 18860
$ws.Range("K131").Value = 2332.7142
$ws.Range("L131").Value = 56580
$ws.Range("M131").Value = 2707.2858
$ws.Range("N131").Value = -66660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3004.6562
$ws.Range("I122").Value = 1948.3334
$ws.Range("J122").Value = 6173.625
$ws.Range("K122").Value = 5845.0002
$ws.Range("L122").Value = 18520.875
$ws.Range("M122").Value = -3395.0002
$ws.Range("N122").Value = -23420.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1418.8462
$ws.Range("I93").Value = 1283.0952
$ws.Range("J93").Value = 1989
$ws.Range("K93").Value = 1283.0952
$ws.Range("L93").Value = 1989
$ws.Range("M93").Value = -35.09519999999998
$ws.Range("N93").Value = -4485
$ws.Range("H120").Value = 93198
$ws.Range("J120").Value = 93198
$ws.Range("L120").Value = 93198
$ws.Range("N120").Value = -102874

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 28025
$ws.Range("J31").Value = 28025
$ws.Range("L31").Value = 28025
$ws.Range("N31").Value = -28721
$ws.Range("H122").Value = 3498.5417
$ws.Range("I122").Value = 3198.476
$ws.Range("J122").Value = 5599
$ws.Range("K122").Value = 9595.428
$ws.Range("L122").Value = 16797
$ws.Range("M122").Value = -7145.428
$ws.Range("N122").Value = -21697
$ws.Range("H136").Value = 8056183
$ws.Range("I136").Value = 10496604
$ws.Range("J136").Value = 2793.8
$ws.Range("K136").Value = 31489812
$ws.Range("L136").Value = 8381.400000000001
$ws.Range("M136").Value = -31487262
$ws.Range("N136").Value = -13481.4
